# Atualizado por script em 01-11-2023 20:45
# Re-orders a few already-scraped fixtures (rows 4/5, 82/83, 94/96, 103/105/106
# get their match data swapped/rotated) and appends two freshly scraped
# fixtures (Brindisi-Catania, Taranto-ACR Messina) as rows 110/111.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 becomes original row 5 content
$row4new = @('Audace Cerignola', 2, 'ACR Messina', 2, 1.68, '01/09/2023 09:12', 1.41, '02/09/2023 20:42', 3.33, '01/09/2023 09:12', 4.29, '02/09/2023 20:44', 4.83, '01/09/2023 09:12', 8.73, '02/09/2023 20:44', 'https://www.betexplorer.com/football/italy/serie-c-group-c/audace-cerignola-acr-messina/h6Hm6UFl/')
for ($i = 0; $i -lt $row4new.Length; $i++) {
    $ws.Cells.Item(4, 6 + $i).Value = $row4new[$i]
}

# Row 5 becomes original row 4 content
$row5new = @('Avellino', 0, 'Latina', 2, 2.03, '01/09/2023 09:12', 1.61, '02/09/2023 20:35', 2.92, '01/09/2023 09:12', 3.7, '02/09/2023 20:35', 3.75, '01/09/2023 09:12', 6.04, '02/09/2023 20:35', 'https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-latina/MuMi5lVf/')
for ($i = 0; $i -lt $row5new.Length; $i++) {
    $ws.Cells.Item(5, 6 + $i).Value = $row5new[$i]
}

# Row 82 becomes original row 83 content
$row82new = @('Giugliano', 0, 'Potenza', 0, 2.79, '19/10/2023 08:13', 2.78, '21/10/2023 18:27', 3.06, '19/10/2023 08:13', 3.38, '21/10/2023 18:23', 2.46, '19/10/2023 08:13', 2.49, '21/10/2023 18:27', 'https://www.betexplorer.com/football/italy/serie-c-group-c/giugliano-potenza/vgGI3jCt/')
for ($i = 0; $i -lt $row82new.Length; $i++) {
    $ws.Cells.Item(82, 6 + $i).Value = $row82new[$i]
}

# Row 83 becomes original row 82 content
$row83new = @('Sorrento', 0, 'Benevento', 1, 3.46, '19/10/2023 08:13', 4.08, '21/10/2023 16:37', 3.07, '19/10/2023 08:13', 3.22, '21/10/2023 16:37', 2.06, '19/10/2023 08:13', 1.99, '21/10/2023 16:37', 'https://www.betexplorer.com/football/italy/serie-c-group-c/sorrento-benevento/hIVN0lda/')
for ($i = 0; $i -lt $row83new.Length; $i++) {
    $ws.Cells.Item(83, 6 + $i).Value = $row83new[$i]
}

# Row 94 becomes original row 96 content
$row94new = @('Casertana', 2, 'Juve Stabia', 1, 2.56, '24/10/2023 12:42', 2.69, '25/10/2023 20:37', 2.83, '24/10/2023 12:42', 3.01, '25/10/2023 20:37', 2.8, '24/10/2023 12:42', 2.83, '25/10/2023 20:37', 'https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-juve-stabia/vL2bC8UL/')
for ($i = 0; $i -lt $row94new.Length; $i++) {
    $ws.Cells.Item(94, 6 + $i).Value = $row94new[$i]
}

# Row 96 becomes original row 94 content
$row96new = @('Taranto', 3, 'Turris', 1, 2.2, '24/10/2023 12:42', 1.79, '25/10/2023 20:44', 3.04, '24/10/2023 12:42', 3.16, '25/10/2023 20:41', 3.24, '24/10/2023 12:42', 5.42, '25/10/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/taranto-turris/vZg5DRu5/')
for ($i = 0; $i -lt $row96new.Length; $i++) {
    $ws.Cells.Item(96, 6 + $i).Value = $row96new[$i]
}

# Row 103 becomes original row 106 content
$row103new = @('Audace Cerignola', 2, 'Casertana', 4, 1.72, '26/10/2023 09:12', 2.19, '29/10/2023 18:21', 3.22, '26/10/2023 09:12', 2.96, '29/10/2023 18:21', 4.74, '26/10/2023 09:12', 3.8, '29/10/2023 18:21', 'https://www.betexplorer.com/football/italy/serie-c-group-c/audace-cerignola-casertana/je62BSqS/')
for ($i = 0; $i -lt $row103new.Length; $i++) {
    $ws.Cells.Item(103, 6 + $i).Value = $row103new[$i]
}

# Row 105 becomes original row 103 content
$row105new = @('Picerno', 2, 'Foggia', 0, 1.98, '27/10/2023 02:42', 3.15, '29/10/2023 18:26', 2.98, '27/10/2023 02:42', 2.75, '29/10/2023 18:26', 3.82, '27/10/2023 02:42', 2.66, '29/10/2023 18:26', 'https://www.betexplorer.com/football/italy/serie-c-group-c/picerno-foggia/nuHWTnQo/')
for ($i = 0; $i -lt $row105new.Length; $i++) {
    $ws.Cells.Item(105, 6 + $i).Value = $row105new[$i]
}

# Row 106 becomes original row 105 content
$row106new = @('Turris', 0, 'Giugliano', 1, 2.13, '27/10/2023 02:42', 1.99, '29/10/2023 18:24', 3.07, '27/10/2023 02:42', 3.54, '29/10/2023 18:24', 3.27, '27/10/2023 02:42', 3.68, '29/10/2023 18:24', 'https://www.betexplorer.com/football/italy/serie-c-group-c/turris-giugliano/K6UwSQfb/')
for ($i = 0; $i -lt $row106new.Length; $i++) {
    $ws.Cells.Item(106, 6 + $i).Value = $row106new[$i]
}

# Insert new row 110 by duplicating formatting from row 109
$ws.Rows.Item(109).Copy()
$ws.Rows.Item(110).Insert()
$ws.Cells.Item(110, 1).Borders.LineStyle = 1
$row110vals = @(109, 'italy', 'serie-c-group-c', '2023-2024', 45231.67708333334, 'Brindisi', 0, 'Catania', 2, 2.79, '07/09/2023 15:12', 3.53, '01/11/2023 15:27', 2.8, '07/09/2023 15:12', 3.07, '01/11/2023 15:27', 2.6, '07/09/2023 15:12', 2.23, '01/11/2023 09:35', 'https://www.betexplorer.com/football/italy/serie-c-group-c/brindisi-catania/vw8NI7af/')
for ($i = 0; $i -lt $row110vals.Length; $i++) {
    $ws.Cells.Item(110, 1 + $i).Value = $row110vals[$i]
}

# Insert new row 111 by duplicating formatting from row 110
$ws.Rows.Item(110).Copy()
$ws.Rows.Item(111).Insert()
$ws.Cells.Item(111, 1).Borders.LineStyle = 1
$row111vals = @(110, 'italy', 'serie-c-group-c', '2023-2024', 45231.77083333334, 'Taranto', 2, 'ACR Messina', 0, 2.07, '15/09/2023 04:42', 1.74, '01/11/2023 17:32', 2.78, '15/09/2023 04:42', 3.34, '01/11/2023 18:20', 4.02, '15/09/2023 04:42', 5.5, '01/11/2023 17:32', 'https://www.betexplorer.com/football/italy/serie-c-group-c/taranto-acr-messina/UmvqePOp/')
for ($i = 0; $i -lt $row111vals.Length; $i++) {
    $ws.Cells.Item(111, 1 + $i).Value = $row111vals[$i]
}
